$d = $word.ActiveDocument

# 1) Merge the runs that were split by proofErr (gramStart/gramEnd, spellStart/spellEnd)
#    markers into single, clean runs. Doing a Find/Replace across the split runs
#    with the full (already-correct) paragraph text forces Word to re-flow the
#    paragraph into one run and drops the now-orphaned proofErr markers.
$mergeTargets = @(
    "On the Map Screen you can press the Scout Button (Which is Rate Limited to once per 6 seconds) in order to look to see if any posts in the area around you match your current hex color.",
    "Hold Down the Scout Button in order to change your color or access a list of Default colors(To make things easier) or topics that you’ve favorited",
    "When you’ve found a post near you (indicated by pins on the map) click it and hit the information button to be sent to the MessageFeed and from there you can read the posts and submit your own. You can even favorite a topic by clicking the star icon at the top."
)
foreach ($t in $mergeTargets) {
    $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null
}

# 2) Merge the "...barebones Post" / " menu." runs (which sandwiched the _GoBack
#    bookmark) into a single run. This also drops the bookmark from this spot,
#    matching the target (a fresh _GoBack bookmark is added further down below,
#    mirroring where the author's edit actually left it).
$postMenuText = "You can submit your own topics as well from the Map Screen just hit the button in the corner and you will be taken to a pretty barebones Post menu."
$d.Content.Find.Execute($postMenuText, $false, $false, $false, $false, $false, $true, 1, $false, $postMenuText, 2) | Out-Null

# 3) Append a new blank paragraph followed by a new closing paragraph about
#    Scout Neon being location based, with a collapsed _GoBack bookmark sitting
#    right before the final "u" (mirroring the original collapsed bookmark that
#    used to sit between "Post" and " menu.").
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' +
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:r><w:t>Scout Neon is location based so will only be able to access posts in the area around yo</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t>u</w:t></w:r>' +
    '</w:p>'
$endRange.InsertXML($newParagraphsXml) | Out-Null
